# New weekly price record for "Locoto" (Agrícola del Norte S.A. de Arica)
# is inserted as row 127, pushing the existing rows 127-150 down to 128-151.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(127).EntireRow.Insert()

$ws.Cells.Item(127, 1).Value  = 1
$ws.Cells.Item(127, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(127, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(127, 4).Value  = 44918
$ws.Cells.Item(127, 5).Value  = 15
$ws.Cells.Item(127, 6).Value  = 100112042
$ws.Cells.Item(127, 7).Value  = "Locoto"
$ws.Cells.Item(127, 8).Value  = "Sin especificar"
$ws.Cells.Item(127, 9).Value  = "Primera"
$ws.Cells.Item(127, 10).Value = 120
$ws.Cells.Item(127, 11).Value = 17000
$ws.Cells.Item(127, 12).Value = 18000
$ws.Cells.Item(127, 13).Value = 17750
$ws.Cells.Item(127, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(127, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(127, 16).Value = 888
$ws.Cells.Item(127, 17).Value = 20
$ws.Cells.Item(127, 18).Value = "Hortaliza"
